$d = $word.ActiveDocument

# Helper: return the 1-based Paragraphs() index whose (trimmed) text
# equals $text exactly. Used so the paragraph-insert steps below don't
# depend on brittle hard-coded positions.
function Get-ParaIndex($text) {
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $t = $d.Paragraphs($i).Range.Text
        $t = $t.TrimEnd([char]13, [char]7)
        if ($t -eq $text) {
            return $i
        }
    }
    return -1
}

# Insert a brand-new "stage direction" paragraph (e.g. "Mom (exit):")
# right after the paragraph whose exact text is $afterText.
function Insert-ParaAfter($afterText, $newText) {
    $idx = Get-ParaIndex $afterText
    $p = $d.Paragraphs($idx)
    $p.Range.InsertParagraphAfter()
    $d.Paragraphs($idx + 1).Range.Text = $newText
}

# --- Simple text substitutions (character-tag swaps / small wording tweaks) ---
# These are all unique substrings so a straightforward Find/Replace is safe.

$d.Content.Find.Execute(
    "Mom (neutral frown): In the end you didn’t end up texting.", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Mom (neutral disappointed): In the end you didn’t end up texting.", 2) | Out-Null

$d.Content.Find.Execute(
    "Mom (neutral smiling): I’m joking, I’m not mad.", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Mom (neutral hehe): I’m joking, I’m not mad.", 2) | Out-Null

$d.Content.Find.Execute(
    "Mom (neutral smiling): Well, they’re basically the same ones.", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Mom (neutral thinking): Well, they’re basically the same ones.", 2) | Out-Null

$d.Content.Find.Execute(
    "Mom (neutral smiling): You should really bring her over for dinner sometime, though.", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Mom (neutral curious): You should really bring her over for dinner sometime, though.", 2) | Out-Null

$d.Content.Find.Execute(
    "Mom (neutral raised_eyebrow): Are you sure?", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Mom (neutral smirk): Are you sure?", 2) | Out-Null

$d.Content.Find.Execute(
    "Mom (neutral smiling): Well, that’s too bad.", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Mom (neutral smiling_eyes_closed): Well, that’s too bad.", 2) | Out-Null

$d.Content.Find.Execute(
    "discuss for so long", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "discuss for what seems like so long", 2) | Out-Null

# --- New stand-alone "stage direction" paragraphs ---
# Inserted from the bottom of the document upward so earlier paragraph
# indices stay valid while each insert is performed.

# After "Mom (neutral smiling): Alright. I'll hold you to that." add "Mom (exit):"
Insert-ParaAfter "Mom (neutral smiling): Alright. I’ll hold you to that." "Mom (exit):"

# Before "I almost choke on my food as well..." (i.e. right after the
# "worried_smile" line) add "Mom (neutral curious):"
Insert-ParaAfter "Mom: (neutral worried_smile) Regardless, you haven’t had a friend over for dinner for a while. It makes me worry sometimes, you know?" "Mom (neutral curious):"

# Before "Pro: Sorry about going out today..." (i.e. right after "I take a
# seat at the kitchen table...") add "Mom (neutral curious):"
Insert-ParaAfter "I take a seat at the kitchen table as my mom heats up dinner in the microwave." "Mom (neutral curious):"

# After "Mom (neutral smiling): I made you dinner anyways..." add "Mom (exit):"
Insert-ParaAfter "Mom (neutral smiling): I made you dinner anyways, so give me a moment to warm it up and then you can eat." "Mom (exit):"
